# The presentation's theme ("Integral") is swapped for the stock "Office
# Theme" palette. In the underlying OOXML this shows up as the full
# contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml trading
# places: the slide master's theme (theme1.xml) becomes the default
# "Office Theme" colours, while the notes master's theme (theme2.xml)
# keeps the "Integral" colours that used to live in theme1.xml.
#
# The PowerPoint object model exposes the live theme's 12 scheme colours
# through Master.Theme.ThemeColorScheme (Item 1..12, in clrScheme document
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Updating them
# here rewrites ppt/theme/theme1.xml's <a:clrScheme> in place.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# dk1 -> #000000 (unchanged)
$tcs.Item(1).RGB = 0
# lt1 -> #FFFFFF (unchanged)
$tcs.Item(2).RGB = 16777215
# dk2 -> #44546A
$tcs.Item(3).RGB = 6968388
# lt2 -> #E7E6E6
$tcs.Item(4).RGB = 15132391
# accent1 -> #5B9BD5
$tcs.Item(5).RGB = 13998939
# accent2 -> #ED7D31
$tcs.Item(6).RGB = 3243501
# accent3 -> #A5A5A5
$tcs.Item(7).RGB = 10855845
# accent4 -> #FFC000
$tcs.Item(8).RGB = 49407
# accent5 -> #4472C4
$tcs.Item(9).RGB = 12874308
# accent6 -> #70AD47
$tcs.Item(10).RGB = 4697456
# hlink -> #0563C1
$tcs.Item(11).RGB = 12673797
# folHlink -> #954F72
$tcs.Item(12).RGB = 7491477
